$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("215")
$ws.Select()

# New row-2 entry for student 21501: name, date (serial 43620 = 2019-06-04), typing score
$ws.Range("B2").Value = "권동훈"

# Reuse the existing date number format from C3 so the new cell shares the same style
# as the other date cells instead of creating a brand-new style entry.
$ws.Range("C3").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("C2").Value = 43620

$ws.Range("D2").Value = 200

# Move the active selection, matching the saved selection in the sheet
$ws.Range("H13").Select()
